$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 857.1429000000001
$ws.Range("I53").Value = 1088.909
$ws.Range("J53").Value = 602.2
$ws.Range("K53").Value = 1088.909
$ws.Range("L53").Value = 602.2
$ws.Range("M53").Value = -451.9090000000001
$ws.Range("N53").Value = -1876.2
$ws.Range("H70").Value = 2379.8
$ws.Range("I70").Value = 1366.3334
$ws.Range("K70").Value = 4099.0002
$ws.Range("M70").Value = -3829.0002
$ws.Range("H73").Value = 2379.8
$ws.Range("I73").Value = 1366.3334
$ws.Range("K73").Value = 4099.0002
$ws.Range("M73").Value = -3163.0002
$ws.Range("H74").Value = 8292.643
$ws.Range("I74").Value = 8084.385
$ws.Range("K74").Value = 8084.385
$ws.Range("M74").Value = -7148.385
$ws.Range("H77").Value = 8292.643
$ws.Range("I77").Value = 8084.385
$ws.Range("K77").Value = 40421.925
$ws.Range("M77").Value = -35741.925
$ws.Range("H125").Value = 1180.5
$ws.Range("I125").Value = 1078
$ws.Range("K125").Value = 9702
$ws.Range("M125").Value = -7242
$ws.Range("H132").Value = 90917020
$ws.Range("I132").Value = 90917020
$ws.Range("K132").Value = 272751060
$ws.Range("M132").Value = -272748530
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 827.25
$ws.Range("I135").Value = 860.63635
$ws.Range("K135").Value = 7745.72715
$ws.Range("M135").Value = -5210.72715
$ws.Range("H137").Value = 12347224
$ws.Range("I137").Value = 20834346
$ws.Range("K137").Value = 62503038
$ws.Range("M137").Value = -62500488
$ws.Range("N134").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1893.7693
$ws.Range("I4").Value = 2582.4285
$ws.Range("J4").Value = 1090.3334
$ws.Range("K4").Value = 2582.4285
$ws.Range("L4").Value = 1090.3334
$ws.Range("M4").Value = -2466.4285
$ws.Range("N4").Value = -1322.3334
$ws.Range("H32").Value = 19695.889
$ws.Range("I32").Value = 22858.084
$ws.Range("J32").Value = 13371.5
$ws.Range("K32").Value = 22858.084
$ws.Range("L32").Value = 13371.5
$ws.Range("M32").Value = -22571.084
$ws.Range("N32").Value = -13945.5
$ws.Range("H54").Value = 14250
$ws.Range("I54").Value = 8500
$ws.Range("K54").Value = 8500
$ws.Range("M54").Value = -7731
$ws.Range("H74").Value = 1553.1904
$ws.Range("I74").Value = 1295.9412
$ws.Range("K74").Value = 1295.9412
$ws.Range("M74").Value = -421.9412
$ws.Range("H77").Value = 1553.1904
$ws.Range("I77").Value = 1295.9412
$ws.Range("K77").Value = 6479.706
$ws.Range("M77").Value = -2111.706
$ws.Range("H132").Value = 4173.125
$ws.Range("I132").Value = 4241
$ws.Range("J132").Value = 2612
$ws.Range("K132").Value = 12723
$ws.Range("L132").Value = 7836
$ws.Range("M132").Value = -10193
$ws.Range("N132").Value = -12896

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2618.1
$ws.Range("I105").Value = 2664.7778
$ws.Range("J105").Value = 2198
$ws.Range("K105").Value = 2664.7778
$ws.Range("L105").Value = 2198
$ws.Range("M105").Value = -917.7777999999998
$ws.Range("N105").Value = -5692
$ws.Range("H134").Value = 1037
$ws.Range("I134").Value = 1046.7858
$ws.Range("K134").Value = 3140.3574
$ws.Range("M134").Value = -605.3574000000003

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 797.3333
$ws.Range("I16").Value = 797.3333
$ws.Range("K16").Value = 797.3333
$ws.Range("M16").Value = -510.3333
$ws.Range("H25").Value = 5048.4443
$ws.Range("I25").Value = 5079.5
$ws.Range("K25").Value = 5079.5
$ws.Range("M25").Value = -4905.5
$ws.Range("H31").Value = 6378.1177
$ws.Range("I31").Value = 5761.933
$ws.Range("K31").Value = 5761.933
$ws.Range("M31").Value = -5466.933
$ws.Range("H34").Value = 6378.1177
$ws.Range("I34").Value = 5761.933
$ws.Range("K34").Value = 5761.933
$ws.Range("M34").Value = -5559.933
$ws.Range("H58").Value = 2231.3914
$ws.Range("I58").Value = 2375
$ws.Range("J58").Value = 1549.25
$ws.Range("K58").Value = 2375
$ws.Range("L58").Value = 1549.25
$ws.Range("M58").Value = -2172
$ws.Range("N58").Value = -1955.25
$ws.Range("H94").Value = 491
$ws.Range("I94").Value = 455
$ws.Range("J94").Value = 599
$ws.Range("K94").Value = 455
$ws.Range("L94").Value = 599
$ws.Range("M94").Value = -4
$ws.Range("N94").Value = -1501
$ws.Range("H99").Value = 3856.0908
$ws.Range("I99").Value = 2764.2
$ws.Range("J99").Value = 4766
$ws.Range("K99").Value = 2764.2
$ws.Range("L99").Value = 4766
$ws.Range("M99").Value = -1266.2
$ws.Range("N99").Value = -7762
$ws.Range("H113").Value = 797.3333
$ws.Range("I113").Value = 797.3333
$ws.Range("K113").Value = 797.3333
$ws.Range("M113").Value = 1372.6667
$ws.Range("H126").Value = 3856.0908
$ws.Range("I126").Value = 2764.2
$ws.Range("J126").Value = 4766
$ws.Range("K126").Value = 8292.599999999999
$ws.Range("L126").Value = 14298
$ws.Range("M126").Value = -5822.599999999999
$ws.Range("N126").Value = -19238
$ws.Range("H134").Value = 2233.2273
$ws.Range("I134").Value = 2407.4
$ws.Range("K134").Value = 7222.200000000001
$ws.Range("M134").Value = -4687.200000000001
$ws.Range("H136").Value = 2231.3914
$ws.Range("I136").Value = 2375
$ws.Range("J136").Value = 1549.25
$ws.Range("K136").Value = 7125
$ws.Range("L136").Value = 4647.75
$ws.Range("M136").Value = -4575
$ws.Range("N136").Value = -9747.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 11547.143
$ws.Range("I56").Value = 11547.143
$ws.Range("K56").Value = 11547.143
$ws.Range("M56").Value = -11017.143
$ws.Range("H122").Value = 932.8889
$ws.Range("I122").Value = 933.3333
$ws.Range("J122").Value = 932
$ws.Range("K122").Value = 8399.9997
$ws.Range("L122").Value = 8388
$ws.Range("M122").Value = -5949.9997
$ws.Range("N122").Value = -13288

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1550.875
$ws.Range("I102").Value = 1589.0952
$ws.Range("K102").Value = 1589.0952
$ws.Range("M102").Value = 32.90480000000002
$ws.Range("H113").Value = 6550.857
$ws.Range("I113").Value = 6766.2
$ws.Range("J113").Value = 6012.5
$ws.Range("K113").Value = 6766.2
$ws.Range("L113").Value = 6012.5
$ws.Range("M113").Value = -4596.2
$ws.Range("N113").Value = -10352.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3501
$ws.Range("J7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("N7").Value = -4224
$ws.Range("H22").Value = 33336040
$ws.Range("I22").Value = 2341.125
$ws.Range("J22").Value = 71431700
$ws.Range("K22").Value = 2341.125
$ws.Range("L22").Value = 71431700
$ws.Range("M22").Value = -2046.125
$ws.Range("N22").Value = -71432290
$ws.Range("H27").Value = 33336040
$ws.Range("I27").Value = 2341.125
$ws.Range("J27").Value = 71431700
$ws.Range("K27").Value = 2341.125
$ws.Range("L27").Value = 71431700
$ws.Range("M27").Value = -2234.125
$ws.Range("N27").Value = -71431914
$ws.Range("H45").Value = 5948
$ws.Range("I45").Value = 3922
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 3922
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -3515
$ws.Range("N45").Value = -10814
$ws.Range("H55").Value = 1378.625
$ws.Range("I55").Value = 1046.3334
$ws.Range("J55").Value = 1578
$ws.Range("K55").Value = 1046.3334
$ws.Range("L55").Value = 1578
$ws.Range("M55").Value = -873.3334
$ws.Range("N55").Value = -1924
$ws.Range("H126").Value = 3501
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 5666.3335
$ws.Range("I132").Value = 5499.5
$ws.Range("K132").Value = 16498.5
$ws.Range("M132").Value = -13968.5
$ws.Range("H136").Value = 498
$ws.Range("I136").Value = 498
$ws.Range("K136").Value = 1494
$ws.Range("M136").Value = 1056

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H81").Value = 6463.706
$ws.Range("I81").Value = 3324
$ws.Range("J81").Value = 13999
$ws.Range("K81").Value = 6648
$ws.Range("L81").Value = 27998
$ws.Range("M81").Value = -5587
$ws.Range("N81").Value = -30120
$ws.Range("H84").Value = 6463.706
$ws.Range("I84").Value = 3324
$ws.Range("J84").Value = 13999
$ws.Range("K84").Value = 33240
$ws.Range("L84").Value = 139990
$ws.Range("M84").Value = -27936
$ws.Range("N84").Value = -150598
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H132").Value = 71445256
$ws.Range("I132").Value = 22464.9
$ws.Range("J132").Value = 250002240
$ws.Range("K132").Value = 67394.70000000001
$ws.Range("L132").Value = 750006720
$ws.Range("M132").Value = -64864.70000000001
$ws.Range("N132").Value = -750011780
$ws.Range("H140").Value = 59396.4
$ws.Range("J140").Value = 59396.4
$ws.Range("L140").Value = 59396.4
$ws.Range("N140").Value = -69756.39999999999
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("N108").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()
